# Append one new record row (row 49) to each of the four data sheets,
# mirroring the prior day's row (row 48) but stamped with the new capture
# timestamp - this is what a new file upload appends to the log.

$wb = $excel.ActiveWorkbook

$newTimestamp = 45835.43721064815

$rows = @(
    @{
        Sheet = "DE_LFT_#1"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x46,0x93,0x3c,0x23,0x3f,0x43,0xe8,0xa0,"
        D = "0x01,0x64"
        E = "0x14"
        F = 380
        G = [double]"7.598631275147109e+23"
        H = 356
        I = 14
    },
    @{
        Sheet = "DE_LFT_#2"
        B = "0x01,0x7c"
        C = "0x00,0xa6,0x60,0x33,0x96,0x39,0x62,0xd0,0x5e,0x78,"
        D = "0x01,0x64"
        E = "0xe"
        F = 380
        G = [double]"5.68432987514711e+23"
        H = 356
        I = 14
    },
    @{
        Sheet = "DE_PLT_#1"
        B = "0x00,0x82"
        C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"
        D = "0x00,0x7F"
        E = "0x7"
        F = 130
        G = [double]"5.68631262647114e+23"
        H = 127
        I = 7
    },
    @{
        Sheet = "DE_PLT_#2"
        B = "0x00,0x82"
        C = "0xd0,0x97,0x78,0x01,0x00,0x00,0x0e,0x3f,0x0c,0x0c,"
        D = "0x00,0x7E"
        E = "0x3"
        F = 130
        G = [double]"9.85046333984776e+23"
        H = 126
        I = 3
    }
)

foreach ($rowData in $rows) {
    $ws = $wb.Worksheets.Item($rowData.Sheet)
    $newRow = 49

    $ws.Cells.Item($newRow, 1).Value = $newTimestamp
    $ws.Cells.Item($newRow, 1).NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $ws.Cells.Item($newRow, 2).Value = $rowData.B
    $ws.Cells.Item($newRow, 3).Value = $rowData.C
    $ws.Cells.Item($newRow, 4).Value = $rowData.D
    $ws.Cells.Item($newRow, 5).Value = $rowData.E
    $ws.Cells.Item($newRow, 6).Value = $rowData.F
    $ws.Cells.Item($newRow, 7).Value = $rowData.G
    $ws.Cells.Item($newRow, 8).Value = $rowData.H
    $ws.Cells.Item($newRow, 9).Value = $rowData.I
}
